$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three coefficient cells to their rounded values
$ws.Range("B3").Value = "-2.82***"
$ws.Range("C2").Value = "-0.01*"
$ws.Range("C3").Value = "-0.47***"
